$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1388.7778
$ws.Range("I28").Value = 999.5
$ws.Range("J28").Value = 1500
$ws.Range("K28").Value = 999.5
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = -514.5
$ws.Range("N28").Value = -2470
$ws.Range("H32").Value = 250000900
$ws.Range("I32").Value = 1000000000
$ws.Range("J32").Value = 1200.6666
$ws.Range("K32").Value = 1000000000
$ws.Range("L32").Value = 1200.6666
$ws.Range("M32").Value = -999999674
$ws.Range("N32").Value = -1852.6666
$ws.Range("H62").Value = 1539.8
$ws.Range("I62").Value = 1674.75
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 1674.75
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = -1050.75
$ws.Range("N62").Value = -2248
$ws.Range("H65").Value = 1539.8
$ws.Range("I65").Value = 1674.75
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 8373.75
$ws.Range("L65").Value = 5000
$ws.Range("M65").Value = -5253.75
$ws.Range("N65").Value = -11240
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H135").Value = 649.8946999999999
$ws.Range("I135").Value = 608.2778
$ws.Range("J135").Value = 1399
$ws.Range("K135").Value = 5474.500199999999
$ws.Range("L135").Value = 12591
$ws.Range("M135").Value = -2939.500199999999
$ws.Range("N135").Value = -17661
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 800
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -501
$ws.Range("H61").Value = 1542.2667
$ws.Range("I61").Value = 1366.7142
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1366.7142
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1154.7142
$ws.Range("N61").Value = -4424
$ws.Range("H74").Value = 2158.814
$ws.Range("I74").Value = 1590.4073
$ws.Range("J74").Value = 3118
$ws.Range("K74").Value = 1590.4073
$ws.Range("L74").Value = 3118
$ws.Range("M74").Value = -716.4073000000001
$ws.Range("N74").Value = -4866
$ws.Range("H77").Value = 2158.814
$ws.Range("I77").Value = 1590.4073
$ws.Range("J77").Value = 3118
$ws.Range("K77").Value = 7952.0365
$ws.Range("L77").Value = 15590
$ws.Range("M77").Value = -3584.0365
$ws.Range("N77").Value = -24326
$ws.Range("H118").Value = 37205
$ws.Range("J118").Value = 37205
$ws.Range("L118").Value = 37205
$ws.Range("N118").Value = -40519
$ws.Range("H125").Value = 42100
$ws.Range("J125").Value = 42100
$ws.Range("L125").Value = 42100
$ws.Range("N125").Value = -51940
$ws.Range("H136").Value = 1542.2667
$ws.Range("I136").Value = 1366.7142
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 4100.142599999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1550.142599999999
$ws.Range("N136").Value = -17100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41565.24
$ws.Range("I20").Value = 54141.633
$ws.Range("J20").Value = 1740
$ws.Range("K20").Value = 54141.633
$ws.Range("L20").Value = 1740
$ws.Range("M20").Value = -53894.633
$ws.Range("N20").Value = -2234
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 4333.3335
$ws.Range("J21").Value = 4333.3335
$ws.Range("L21").Value = 4333.3335
$ws.Range("N21").Value = -4803.3335
$ws.Range("H58").Value = 10187.896
$ws.Range("J58").Value = 19415.285
$ws.Range("L58").Value = 19415.285
$ws.Range("N58").Value = -19821.285
$ws.Range("H132").Value = 4361.92
$ws.Range("I132").Value = 5214
$ws.Range("K132").Value = 15642
$ws.Range("M132").Value = -13112
$ws.Range("H134").Value = 2097.2856
$ws.Range("I134").Value = 2236.2
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 6708.599999999999
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -4173.599999999999
$ws.Range("N134").Value = -10320
$ws.Range("H136").Value = 10187.896
$ws.Range("J136").Value = 19415.285
$ws.Range("L136").Value = 58245.855
$ws.Range("N136").Value = -63345.855
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 29986.857
$ws.Range("I5").Value = 1977
$ws.Range("J5").Value = 67333.336
$ws.Range("K5").Value = 5931
$ws.Range("L5").Value = 202000.008
$ws.Range("M5").Value = -5819
$ws.Range("N5").Value = -202224.008
$ws.Range("H122").Value = 7801.2144
$ws.Range("J122").Value = 21079.6
$ws.Range("L122").Value = 189716.4
$ws.Range("N122").Value = -194616.4
$ws.Range("H131").Value = 848.4400000000001
$ws.Range("I131").Value = 484.2857
$ws.Range("J131").Value = 875.8495
$ws.Range("K131").Value = 1452.8571
$ws.Range("L131").Value = 2627.5485
$ws.Range("M131").Value = 3587.1429
$ws.Range("N131").Value = -12707.5485
$ws.Range("H135").Value = 29986.857
$ws.Range("I135").Value = 1977
$ws.Range("J135").Value = 67333.336
$ws.Range("K135").Value = 17793
$ws.Range("L135").Value = 606000.024
$ws.Range("M135").Value = -15258
$ws.Range("N135").Value = -611070.024
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2157.7827
$ws.Range("I126").Value = 1939.4286
$ws.Range("J126").Value = 2253.3125
$ws.Range("K126").Value = 5818.2858
$ws.Range("L126").Value = 6759.9375
$ws.Range("M126").Value = -3348.2858
$ws.Range("N126").Value = -11699.9375
$ws.Range("H132").Value = 3209.8096
$ws.Range("I132").Value = 2725.7273
$ws.Range("K132").Value = 8177.1819
$ws.Range("M132").Value = -5647.1819
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5651.6924
$ws.Range("I132").Value = 6402.6665
$ws.Range("K132").Value = 19207.9995
$ws.Range("M132").Value = -16677.9995
$ws.Range("H136").Value = 1488.3889
$ws.Range("I136").Value = 1472.8182
$ws.Range("K136").Value = 4418.4546
$ws.Range("M136").Value = -1868.4546
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 48892.5
$ws.Range("J116").Value = 48892.5
$ws.Range("L116").Value = 48892.5
$ws.Range("N116").Value = -58070.5
$ws.Range("H117").Value = 39900
$ws.Range("J117").Value = 39900
$ws.Range("L117").Value = 39900
$ws.Range("N117").Value = -49078
$ws.Range("H119").Value = 41979.668
$ws.Range("J119").Value = 41979.668
$ws.Range("L119").Value = 41979.668
$ws.Range("N119").Value = -51655.668
$ws.Range("H120").Value = 36990
$ws.Range("J120").Value = 36990
$ws.Range("L120").Value = 36990
$ws.Range("N120").Value = -46666
$ws.Range("H125").Value = 39949.5
$ws.Range("J125").Value = 39949.5
$ws.Range("L125").Value = 39949.5
$ws.Range("N125").Value = -49789.5
$ws.Range("H136").Value = 1272.3024
$ws.Range("I136").Value = 500.4
$ws.Range("J136").Value = 2344.389
$ws.Range("K136").Value = 1501.2
$ws.Range("L136").Value = 7033.167
$ws.Range("M136").Value = 1048.8
$ws.Range("N136").Value = -12133.167
